$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "305.28"
Set-TextValue $ws.Range("E2") "-0.24%"
Set-TextValue $ws.Range("G2") "9"

Set-TextValue $ws.Range("E3") "-0.86%"
Set-TextValue $ws.Range("G3") "9"

Set-TextValue $ws.Range("D4") "5.049"
Set-TextValue $ws.Range("E4") "-0.95%"
Set-TextValue $ws.Range("G4") "9"

Set-TextValue $ws.Range("D5") "0.08036"
Set-TextValue $ws.Range("E5") "-0.70%"
Set-TextValue $ws.Range("G5") "9"

Set-TextValue $ws.Range("E6") "-1.37%"
Set-TextValue $ws.Range("G6") "9"

Set-TextValue $ws.Range("B7") "GateToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "4.157"
Set-TextValue $ws.Range("E7") "-0.75%"
Set-TextValue $ws.Range("G7") "9"

Set-TextValue $ws.Range("B8") "KuCoinToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D8") "7.780"
Set-TextValue $ws.Range("E8") "0.31%"
Set-TextValue $ws.Range("G8") "9"

Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9215"
Set-TextValue $ws.Range("E9") "-0.77%"
Set-TextValue $ws.Range("G9") "9"

Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1278"
Set-TextValue $ws.Range("E10") "-6.61%"
Set-TextValue $ws.Range("G10") "9"

Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1915"
Set-TextValue $ws.Range("E11") "0.44%"
Set-TextValue $ws.Range("G11") "9"

Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09057"
Set-TextValue $ws.Range("E12") "-1.57%"
Set-TextValue $ws.Range("G12") "9"

Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03456"
Set-TextValue $ws.Range("E13") "0.93%"
Set-TextValue $ws.Range("G13") "9"

Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09853"
Set-TextValue $ws.Range("E14") "0.20%"
Set-TextValue $ws.Range("G14") "9"

Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001412"
Set-TextValue $ws.Range("E15") "-0.79%"
Set-TextValue $ws.Range("G15") "9"

Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006206"
Set-TextValue $ws.Range("E16") "6.28%"
Set-TextValue $ws.Range("G16") "9"

Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.736"
Set-TextValue $ws.Range("E17") "3.58%"
Set-TextValue $ws.Range("G17") "9"

Set-TextValue $ws.Range("D18") "3.395"
Set-TextValue $ws.Range("E18") "12.61%"
Set-TextValue $ws.Range("G18") "9"

Set-TextValue $ws.Range("D19") "0.3446"
Set-TextValue $ws.Range("E19") "-0.26%"
Set-TextValue $ws.Range("G19") "9"

Set-TextValue $ws.Range("D20") "0.1349"
Set-TextValue $ws.Range("E20") "1.20%"
Set-TextValue $ws.Range("G20") "9"

Set-TextValue $ws.Range("D21") "5.169"
Set-TextValue $ws.Range("E21") "5.53%"
Set-TextValue $ws.Range("G21") "9"

Set-TextValue $ws.Range("D22") "0.2605"
Set-TextValue $ws.Range("E22") "1.46%"
Set-TextValue $ws.Range("G22") "9"

Set-TextValue $ws.Range("D23") "0.04432"
Set-TextValue $ws.Range("E23") "-0.02%"
Set-TextValue $ws.Range("G23") "9"

Set-TextValue $ws.Range("E24") "1.08%"
Set-TextValue $ws.Range("G24") "9"

Set-TextValue $ws.Range("D25") "0.004616"
Set-TextValue $ws.Range("E25") "-4.15%"
Set-TextValue $ws.Range("G25") "9"

Set-TextValue $ws.Range("G26") "9"

Set-TextValue $ws.Range("D27") "0.0001253"
Set-TextValue $ws.Range("E27") "-3.65%"
Set-TextValue $ws.Range("G27") "9"

Set-TextValue $ws.Range("D28") "0.0004450"
Set-TextValue $ws.Range("E28") "42.02%"
Set-TextValue $ws.Range("G28") "9"

Set-TextValue $ws.Range("G29") "9"

Set-TextValue $ws.Range("G30") "9"

Set-TextValue $ws.Range("G31") "9"

Set-TextValue $ws.Range("G32") "9"

Set-TextValue $ws.Range("G33") "9"

Set-TextValue $ws.Range("G34") "9"

Set-TextValue $ws.Range("G35") "9"

Set-TextValue $ws.Range("G36") "9"

Set-TextValue $ws.Range("G37") "9"

Set-TextValue $ws.Range("G38") "9"

Set-TextValue $ws.Range("D39") "0.01945"
Set-TextValue $ws.Range("E39") "-3.66%"
Set-TextValue $ws.Range("G39") "9"

Set-TextValue $ws.Range("D40") "0.05661"
Set-TextValue $ws.Range("E40") "14.97%"
Set-TextValue $ws.Range("G40") "9"

Set-TextValue $ws.Range("D41") "0.007627"
Set-TextValue $ws.Range("E41") "0.17%"
Set-TextValue $ws.Range("G41") "9"

Set-TextValue $ws.Range("D42") "0.01017"
Set-TextValue $ws.Range("E42") "-2.12%"
Set-TextValue $ws.Range("G42") "9"

Set-TextValue $ws.Range("E43") "-1.68%"
Set-TextValue $ws.Range("G43") "9"

Set-TextValue $ws.Range("D44") "0.002175"
Set-TextValue $ws.Range("E44") "3.43%"
Set-TextValue $ws.Range("G44") "9"

Set-TextValue $ws.Range("D45") "0.009844"
Set-TextValue $ws.Range("E45") "-10.69%"
Set-TextValue $ws.Range("G45") "9"

Set-TextValue $ws.Range("D46") "0.00006147"
Set-TextValue $ws.Range("E46") "-4.28%"
Set-TextValue $ws.Range("G46") "9"

Set-TextValue $ws.Range("E47") "0.03%"
Set-TextValue $ws.Range("G47") "9"

Set-TextValue $ws.Range("D48") "63.57"
Set-TextValue $ws.Range("E48") "-1.69%"
Set-TextValue $ws.Range("G48") "9"

Set-TextValue $ws.Range("D49") "0.001661"
Set-TextValue $ws.Range("E49") "39.34%"
Set-TextValue $ws.Range("G49") "9"

Set-TextValue $ws.Range("D50") "0.00002104"
Set-TextValue $ws.Range("E50") "0.03%"
Set-TextValue $ws.Range("G50") "9"

Set-TextValue $ws.Range("D51") "0.0002004"
Set-TextValue $ws.Range("E51") "0.03%"
Set-TextValue $ws.Range("G51") "9"
